# Insert a new data row at row 165 (pushing the existing rows 165:204 down to
# 166:205) and populate it with the new price-observation record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(165).Insert()

$ws.Cells.Item(165, 1).Value = 7
$ws.Cells.Item(165, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(165, 3).Value = "Ñuble"
$ws.Cells.Item(165, 4).Value = 44641
$ws.Cells.Item(165, 5).Value = 16
$ws.Cells.Item(165, 6).Value = 100112032
$ws.Cells.Item(165, 7).Value = "Zapallo italiano"
$ws.Cells.Item(165, 8).Value = "Sin especificar"
$ws.Cells.Item(165, 9).Value = "Primera"
$ws.Cells.Item(165, 10).Value = 120
$ws.Cells.Item(165, 11).Value = 7500
$ws.Cells.Item(165, 12).Value = 8000
$ws.Cells.Item(165, 13).Value = 7750
$ws.Cells.Item(165, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(165, 15).Value = "Región del Maule"
$ws.Cells.Item(165, 16).Value = 155
$ws.Cells.Item(165, 17).Value = 50
$ws.Cells.Item(165, 18).Value = "Hortaliza"
